$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.771.97'
$ws.Range("E2").Value = '  -6.95%  '

$ws.Range("D3").Value = '2.533.33'
$ws.Range("E3").Value = '  -3.26%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '298.63'
$ws.Range("E5").Value = '  -3.47%  '

$ws.Range("D6").Value = '91.74'
$ws.Range("E6").Value = '  -6.67%  '

$ws.Range("E7").Value = '  -3.72%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").Value = '0.549'
$ws.Range("E9").Value = '  -5.12%  '

$ws.Range("D10").Value = '35.97'
$ws.Range("E10").Value = '  -6.98%  '

$ws.Range("D11").Value = '0.0803'
$ws.Range("E11").Value = '  -4.54%  '

$ws.Range("D12").Value = '7.65'
$ws.Range("E12").Value = '  -5.10%  '

$ws.Range("E13").Value = '  +6.00%  '

$ws.Range("D14").Value = '2.920.14'
$ws.Range("E14").Value = '  -3.22%  '

$ws.Range("D15").Value = '2.531.59'
$ws.Range("E15").Value = '  -3.13%  '

$ws.Range("E16").Value = '  -5.23%  '

$ws.Range("D17").Value = '14.07'
$ws.Range("E17").Value = '  -5.00%  '

$ws.Range("D18").Value = '42.814.96'
$ws.Range("E18").Value = '  -6.82%  '

$ws.Range("D19").Value = '13.03'
$ws.Range("E19").Value = '  +2.49%  '

$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").Value = '  -3.04%  '

$ws.Range("D21").Value = '6.53'
$ws.Range("E21").Value = '  -3.32%  '

$ws.Range("D22").Value = '71.42'
$ws.Range("E22").Value = '  -3.99%  '

$ws.Range("D23").Value = '256.30'
$ws.Range("E23").Value = '  -9.34%  '

$ws.Range("D24").Value = '2.90'
$ws.Range("E24").Value = '  -4.15%  '

$ws.Range("D25").Value = '29.32'
$ws.Range("E25").Value = '  -0.80%  '

$ws.Range("D26").Value = '2.10'
$ws.Range("E26").Value = '  -6.73%  '

$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -4.77%  '

$ws.Range("D29").Value = '36.89'
$ws.Range("E29").Value = '  -4.40%  '

$ws.Range("E30").Value = '  -3.79%  '

$ws.Range("E31").Value = '  -4.67%  '

$ws.Range("D32").Value = '152.13'
$ws.Range("E32").Value = '  -2.91%  '

$ws.Range("E33").Value = '  -5.05%  '

$ws.Range("E34").Value = '  -2.00%  '

$ws.Range("E35").Value = '  -6.19%  '

$ws.Range("D36").Value = '0.0791'
$ws.Range("E36").Value = '  -5.44%  '

$ws.Range("E37").Value = '  -6.73%  '

$ws.Range("D38").Value = '0.119'
$ws.Range("E38").Value = '  -3.30%  '

$ws.Range("D39").Value = '23.95'
$ws.Range("E39").Value = '  +7.66%  '

$ws.Range("D40").Value = '16.78'
$ws.Range("E40").Value = '  +6.08%  '

$ws.Range("D41").Value = '3.88'
$ws.Range("E41").Value = '  -3.79%  '

$ws.Range("D42").Value = '0.0309'
$ws.Range("E42").Value = '  -5.19%  '

$ws.Range("D43").Value = '3.37'
$ws.Range("E43").Value = '  -4.57%  '

$ws.Range("D44").Value = '2.086.46'
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = '84.05'
$ws.Range("E46").Value = '  -10.48%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '8.97'
$ws.Range("E47").Value = '  -1.82%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '1.60'
$ws.Range("E48").Value = '  +1.72%  '

$ws.Range("D49").Value = '2.778.07'
$ws.Range("E49").Value = '  -3.18%  '

$ws.Range("D50").Value = '103.92'
$ws.Range("E50").Value = '  -5.54%  '

$ws.Range("E51").Value = '  -4.62%  '
